$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value (correct marks per question)
$ws.Range("B11").Value = 5

# Update total correct marks
$ws.Range("B12").Value = 70

# Update "Corr/total" marks display
$ws.Range("E12").Value = "70/140"
